$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Merge the "Example 4:  Confidence Interval ..." heading runs into one.
#    The visible text does not change, only the run structure (5 runs -> 1).
#    Re-running Find/Replace over the identical text collapses the runs into
#    a single run while keeping the formatting (rPr) of the first run.
# ---------------------------------------------------------------------------
$titleText = "Example 4:  Confidence Interval for the Difference of Two Sample Proportions"
$range = $d.Content
$range.Find.Execute($titleText, $false, $false, $false, $false, $false, $true, 0, $false, $titleText, 1) | Out-Null

# ---------------------------------------------------------------------------
# 2) Strip the underscores from the R variable names in the source-code
#    blocks (x_1 -> x1, n_1 -> n1, phat_1 -> phat1, x_2 -> x2, n_2 -> n2,
#    phat_2 -> phat2), walking strictly left-to-right through the document
#    so each Find/Replace call lands on the intended occurrence.
# ---------------------------------------------------------------------------
$range = $d.Content

$replacements = @(
    @("x_1 ", "x1 "),
    @("n_1 ", "n1 "),
    @("phat_1 ", "phat1 "),
    @(" x_1 ", " x1 "),
    @(" n_1", " n1"),
    @("x_2 ", "x2 "),
    @("n_2 ", "n2 "),
    @("phat_2 ", "phat2 "),
    @(" x_2 ", " x2 "),
    @(" n_2", " n2"),
    @(" phat_1 ", " phat1 "),
    @(" phat_2", " phat2"),
    @("((phat_1 ", "((phat1 "),
    @(" phat_1) ", " phat1) ")
)

foreach ($pair in $replacements) {
    $range.Find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 0, $false, $pair[1], 1) | Out-Null
}

# ---------------------------------------------------------------------------
# 3) The next token " n_2) " is a special case: the source actually had a
#    typo (n_2 where n_1 was meant), so the edit turns it into " n1) " while
#    also splitting the previously single run into three runs (" n", "1",
#    ") "), each still styled NormalTok. Delete the "_2" and retype "1",
#    nudging Bold on/off on just the new character so the run isn't folded
#    back into its neighbours.
# ---------------------------------------------------------------------------
$range.Find.Execute("_2) ", $false, $false, $false, $false, $false, $true, 0, $false, ") ", 1) | Out-Null
$digitStart = $range.Start
$digitRange = $d.Range($digitStart, $digitStart)
$digitRange.InsertBefore("1")
$digitRange = $d.Range($digitStart, $digitStart + 1)
$digitRange.Bold = 1
$digitRange.Bold = 0

# ---------------------------------------------------------------------------
# 4) Continue the underscore clean-up for the remaining formula tokens.
# ---------------------------------------------------------------------------
$range = $d.Range($digitStart + 1, $d.Content.End)

$replacements2 = @(
    @(" (phat_2 ", " (phat2 "),
    @(" phat_2) ", " phat2) "),
    @(" n_2))", " n2))")
)

foreach ($pair in $replacements2) {
    $range.Find.Execute($pair[0], $false, $false, $false, $false, $false, $true, 0, $false, $pair[1], 1) | Out-Null
}
